$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sperimentazioni 5 vicini")
$rng = $ws.Range("A27:Q27")
$rng.Interior.TintAndShade = -0.34998626667073579
$rng.Interior.ThemeColor = 2
